# Update "想去人数" (interest/attendance count) figures that changed between
# the previous data pull and the latest generated output (commit 456a3b4).
#
# Sheet "展览" (Exhibition)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 132
$ws1.Range("F8").Value  = 3621
$ws1.Range("F10").Value = 2473
$ws1.Range("F12").Value = 2898
$ws1.Range("F15").Value = 2227
$ws1.Range("F17").Value = 105
$ws1.Range("F21").Value = 168
$ws1.Range("F24").Value = 267
$ws1.Range("F28").Value = 1273
$ws1.Range("F32").Value = 3977
$ws1.Range("F33").Value = 3451
$ws1.Range("F34").Value = 46
$ws1.Range("F37").Value = 418
$ws1.Range("F45").Value = 27

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 2143

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 132
$ws4.Range("F14").Value = 3621
$ws4.Range("F15").Value = 2473
$ws4.Range("F17").Value = 2898
$ws4.Range("F19").Value = 2227
$ws4.Range("F21").Value = 105
$ws4.Range("F27").Value = 267
$ws4.Range("F31").Value = 1273
$ws4.Range("F36").Value = 3977
$ws4.Range("F37").Value = 3451
$ws4.Range("F38").Value = 46
$ws4.Range("F41").Value = 418
$ws4.Range("F48").Value = 27
